$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 54 (pushes old row 54 down to become row 56)
$ws.Rows("54:55").Insert()

# Update shifted/changed values for rows 33-53
# Row 33
$ws.Range("D33").Value = 44554
$ws.Range("M33").Value = 50
$ws.Range("N33").Value = 6000
$ws.Range("O33").Value = 6000
$ws.Range("P33").Value = 6000
$ws.Range("S33").Value = 3000

# Row 34
$ws.Range("D34").Value = 44181
$ws.Range("M34").Value = 140
$ws.Range("N34").Value = 4000
$ws.Range("O34").Value = 4500
$ws.Range("P34").Value = 4250
$ws.Range("S34").Value = 2125

# Row 35
$ws.Range("D35").Value = 44209
$ws.Range("M35").Value = 170
$ws.Range("N35").Value = 3000
$ws.Range("P35").Value = 3500
$ws.Range("S35").Value = 1750

# Row 36
$ws.Range("D36").Value = 44176
$ws.Range("M36").Value = 100
$ws.Range("N36").Value = 4000
$ws.Range("O36").Value = 4000
$ws.Range("P36").Value = 4000
$ws.Range("R36").Value = 'Provincia de Linares'
$ws.Range("S36").Value = 2000

# Row 37
$ws.Range("D37").Value = 44539
$ws.Range("M37").Value = 500
$ws.Range("N37").Value = 5000
$ws.Range("O37").Value = 5000
$ws.Range("P37").Value = 5000
$ws.Range("R37").Value = 'Provincia de Curicó'
$ws.Range("S37").Value = 2500

# Row 38
$ws.Range("D38").Value = 44210
$ws.Range("M38").Value = 400
$ws.Range("N38").Value = 3000
$ws.Range("P38").Value = 3500
$ws.Range("R38").Value = 'Provincia de Linares'
$ws.Range("S38").Value = 1750

# Row 39
$ws.Range("D39").Value = 44250
$ws.Range("M39").Value = 100

# Row 40
$ws.Range("D40").Value = 44553
$ws.Range("M40").Value = 250
$ws.Range("R40").Value = 'Provincia de Curicó'

# Row 41
$ws.Range("D41").Value = 44251
$ws.Range("M41").Value = 125
$ws.Range("N41").Value = 4000
$ws.Range("O41").Value = 4000
$ws.Range("P41").Value = 4000
$ws.Range("S41").Value = 2000

# Row 42
$ws.Range("D42").Value = 44551
$ws.Range("M42").Value = 500
$ws.Range("N42").Value = 6000
$ws.Range("O42").Value = 6000
$ws.Range("P42").Value = 6000
$ws.Range("S42").Value = 3000

# Row 43
$ws.Range("L43").Value = 'Primera'
$ws.Range("M43").Value = 100
$ws.Range("N43").Value = 3400
$ws.Range("O43").Value = 3400
$ws.Range("P43").Value = 3400
$ws.Range("R43").Value = 'Provincia de Curicó'
$ws.Range("S43").Value = 1700

# Row 44
$ws.Range("D44").Value = 44187
$ws.Range("M44").Value = 200
$ws.Range("R44").Value = 'Provincia de Linares'

# Row 45
$ws.Range("D45").Value = 44187
$ws.Range("L45").Value = 'Segunda'
$ws.Range("M45").Value = 50
$ws.Range("N45").Value = 3000
$ws.Range("O45").Value = 3000
$ws.Range("P45").Value = 3000
$ws.Range("S45").Value = 1500

# Row 46
$ws.Range("D46").Value = 44257
$ws.Range("M46").Value = 100
$ws.Range("R46").Value = 'Provincia de Curicó'

# Row 47
$ws.Range("D47").Value = 44175
$ws.Range("M47").Value = 250
$ws.Range("R47").Value = 'Provincia de Linares'

# Row 48
$ws.Range("D48").Value = 44196
$ws.Range("M48").Value = 150

# Row 49
$ws.Range("D49").Value = 44188
$ws.Range("M49").Value = 300

# Row 50
$ws.Range("D50").Value = 44188
$ws.Range("M50").Value = 500

# Row 51
$ws.Range("D51").Value = 44224
$ws.Range("M51").Value = 250
$ws.Range("N51").Value = 4000
$ws.Range("O51").Value = 4000
$ws.Range("P51").Value = 4000
$ws.Range("R51").Value = 'Provincia de Curicó'
$ws.Range("S51").Value = 2000

# Row 52
$ws.Range("D52").Value = 44224
$ws.Range("M52").Value = 300
$ws.Range("N52").Value = 4000
$ws.Range("P52").Value = 4000
$ws.Range("R52").Value = 'Provincia de Linares'
$ws.Range("S52").Value = 2000

# Row 53
$ws.Range("D53").Value = 44195
$ws.Range("M53").Value = 300
$ws.Range("N53").Value = 3000
$ws.Range("O53").Value = 3000
$ws.Range("P53").Value = 3000
$ws.Range("S53").Value = 1500

# Fully populate new rows 54 and 55 (old row 54 automatically shifted to row 56 by Insert above)
# Row 54
$ws.Range("A54").Value = 6
$ws.Range("B54").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C54").Value = 'Metropolitana'
$ws.Range("D54").Value = 44239
$ws.Range("E54").Value = 13
$ws.Range("F54").Value = 'Fruta'
$ws.Range("G54").Value = 100101
$ws.Range("H54").Value = 'Berries'
$ws.Range("I54").Value = 100101008
$ws.Range("J54").Value = 'Mora'
$ws.Range("K54").Value = 'Sin especificar'
$ws.Range("L54").Value = 'Primera'
$ws.Range("M54").Value = 350
$ws.Range("N54").Value = 3500
$ws.Range("O54").Value = 4000
$ws.Range("P54").Value = 3750
$ws.Range("Q54").Value = '$/bandeja 2 kilos'
$ws.Range("R54").Value = 'Provincia de Curicó'
$ws.Range("S54").Value = 1875
$ws.Range("T54").Value = 2

# Row 55
$ws.Range("A55").Value = 6
$ws.Range("B55").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C55").Value = 'Metropolitana'
$ws.Range("D55").Value = 44194
$ws.Range("E55").Value = 13
$ws.Range("F55").Value = 'Fruta'
$ws.Range("G55").Value = 100101
$ws.Range("H55").Value = 'Berries'
$ws.Range("I55").Value = 100101008
$ws.Range("J55").Value = 'Mora'
$ws.Range("K55").Value = 'Sin especificar'
$ws.Range("L55").Value = 'Primera'
$ws.Range("M55").Value = 250
$ws.Range("N55").Value = 4000
$ws.Range("O55").Value = 4000
$ws.Range("P55").Value = 4000
$ws.Range("Q55").Value = '$/bandeja 2 kilos'
$ws.Range("R55").Value = 'Provincia de Linares'
$ws.Range("S55").Value = 2000
$ws.Range("T55").Value = 2
